# added data retrieval from api request to eleven travel
# -----------------------------------------------------------------
# The "prijs" column (column D) held price strings scraped in the
# old format, e.g. "€ 44,00" (a euro sign followed by a non-breaking
# space and the amount). The new API-driven retrieval stores just
# the bare amount, e.g. "44,00". Strip the leading currency prefix
# from every price cell in column D (row 1 is the header "prijs"
# and is left alone).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$euro = [char]0x20AC
$nbsp = [char]0xA0
$prefixNbsp = "$euro$nbsp"
$prefixSpace = "$euro "

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $v = $cell.Value2

    if ($v -eq $null) {
        continue
    }

    $text = [string]$v

    if ($text.StartsWith($prefixNbsp)) {
        $cell.Value = $text.Substring($prefixNbsp.Length)
    } elseif ($text.StartsWith($prefixSpace)) {
        $cell.Value = $text.Substring($prefixSpace.Length)
    } elseif ($text.StartsWith($euro)) {
        $cell.Value = $text.Substring(1).TrimStart($nbsp).TrimStart(" ")
    }
}
